$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set header values for the two new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting of the existing header cell (H1) onto the new headers
# so they get the same bold/centered/bordered style (cellXf index 1).
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("J1").PasteSpecial(-4122)

# Add the data values for rows 2 and 3 in the two new columns
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9

$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 8

$excel.CutCopyMode = 0
